$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.879.17"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "2.440.39"
$ws.Range("E3").Value = "  -1.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.27"
$ws.Range("E5").Value = "  -0.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.13"
$ws.Range("E6").Value = "  -1.08%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.47%  "

# Row 9
$ws.Range("E9").Value = "  +8.13%  "

# Row 10
$ws.Range("E10").Value = "  -2.08%  "

# Row 11
$ws.Range("E11").Value = "  -0.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.60"
$ws.Range("E12").Value = "  -5.22%  "

# Row 13
$ws.Range("E13").Value = "  +3.94%  "

# Row 14
$ws.Range("D14").Value = "68.759.71"
$ws.Range("E14").Value = "  +0.08%  "

# Row 15
$ws.Range("D15").Value = "2.887.18"
$ws.Range("E15").Value = "  -0.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.27"
$ws.Range("E16").Value = "  -1.45%  "

# Row 17
$ws.Range("D17").Value = "2.437.70"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.59"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.58"
$ws.Range("E19").Value = "  +0.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.95"
$ws.Range("E20").Value = "  +0.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.84"
$ws.Range("E21").Value = "  +1.05%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.93"
$ws.Range("E22").Value = "  +2.42%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.09"
$ws.Range("E24").Value = "  +0.60%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.71"
$ws.Range("E25").Value = "  +1.09%  "

# Row 26
$ws.Range("D26").Value = "2.564.31"
$ws.Range("E26").Value = "  -1.27%  "

# Row 27
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.23"
$ws.Range("E27").Value = "  -0.31%  "

# Row 28
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.953"
$ws.Range("E28").Value = "  -4.13%  "

# Row 29
$ws.Range("E29").Value = "  -0.58%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.14"
$ws.Range("E30").Value = "  -1.28%  "

# Row 31
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  +1.38%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "428.70"
$ws.Range("E33").Value = "  -0.45%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.60"
$ws.Range("E34").Value = "  -1.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.81"
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.99"
$ws.Range("E36").Value = "  -0.12%  "

# Row 37
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.03"
$ws.Range("E38").Value = "  +0.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.106"
$ws.Range("E39").Value = "  -2.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.300"
$ws.Range("E40").Value = "  +0.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.51"
$ws.Range("E41").Value = "  +2.18%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.35"
$ws.Range("E42").Value = "  -2.40%  "

# Row 43
$ws.Range("E43").Value = "  -1.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.06"
$ws.Range("E44").Value = "  -0.46%  "

# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.33"
$ws.Range("E45").Value = "  -1.81%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.32"
$ws.Range("E46").Value = "  -0.11%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0721"
$ws.Range("E47").Value = "  +0.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.481"
$ws.Range("E48").Value = "  -1.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.559"
$ws.Range("E49").Value = "  -1.04%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0923"
$ws.Range("E50").Value = "  +0.47%  "

# Row 51
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.15"
$ws.Range("E51").Value = "  +2.81%  "
